{"js": "// Kien's report: bump the development+testing cost from 12.000.000 VND to\n// 20.000.000 VND in the \"\u01af\u1edbc l\u01b0\u1ee3ng gi\u00e1 th\u00e0nh\" section.\n//\n// Strategy: locate the unique sentence, then narrow the search to the\n// \"12\" that prefixes the amount and replace just that substring. This\n// keeps every other character (and the run's Times New Roman / 28pt\n// formatting) untouched, mirroring the author's \"12\" -> \"20\" edit.\n\nconst sentence = \"Chi ph\u00ed ph\u00e1t tri\u1ec3n  +  Chi ph\u00ed ki\u1ec3m th\u1eed: 12.000.000 VND\";\n\nconst matches = context.document.body.search(sentence, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find the cost sentence to update.\");\n}\n\nconst target = matches.items[0];\n\nconst amountMatches = target.search(\"12.000.000 VND\", { matchCase: true });\namountMatches.load(\"items\");\nawait context.sync();\n\nif (amountMatches.items.length === 0) {\n  throw new Error(\"Could not find the '12.000.000 VND' amount to update.\");\n}\n\nconst amountRange = amountMatches.items[0];\n\nconst prefixMatches = amountRange.search(\"12\", { matchCase: true });\nprefixMatches.load(\"items\");\nawait context.sync();\n\nif (prefixMatches.items.length === 0) {\n  throw new Error(\"Could not find the '12' prefix of the amount to update.\");\n}\n\nprefixMatches.items[0].insertText(\"20\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Kien's report: bump the development+testing cost from 12.000.000 VND to\n# 20.000.000 VND in the \"Uoc luong gia thanh\" section.\n#\n# Strategy: locate the unique sentence via Find, then narrow a duplicated\n# range to just the \"12\" that prefixes the amount and replace that\n# substring's text. This keeps every other character (and the run's\n# Times New Roman / 28pt formatting) untouched, mirroring the author's\n# \"12\" -> \"20\" edit.\n\n$d = $word.ActiveDocument\n\n$sentence = \"Chi ph\u00ed ph\u00e1t tri\u1ec3n  +  Chi ph\u00ed ki\u1ec3m th\u1eed: 12.000.000 VND\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute($sentence)\n\nif (-not $found) {\n    throw \"Could not find the cost sentence to update.\"\n}\n\n# $rng now spans the matched sentence; narrow a copy to the amount.\n$amountRng = $rng.Duplicate\n$amountFind = $amountRng.Find\n$amountFind.MatchCase = $true\n$amountFound = $amountFind.Execute(\"12.000.000 VND\")\n\nif (-not $amountFound) {\n    throw \"Could not find the '12.000.000 VND' amount to update.\"\n}\n\n# $amountRng now spans \"12.000.000 VND\"; narrow further to the \"12\" prefix.\n$prefixRng = $amountRng.Duplicate\n$prefixFind = $prefixRng.Find\n$prefixFind.MatchCase = $true\n$prefixFound = $prefixFind.Execute(\"12\")\n\nif (-not $prefixFound) {\n    throw \"Could not find the '12' prefix of the amount to update.\"\n}\n\n$prefixRng.Text = \"20\"\n"}
